$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 3-13 hold the augment data. Columns D and F drop from 10 to 0,
# and column G rises from 10 to 15. Column E is left unchanged at 10.
for ($r = 3; $r -le 13; $r++) {
    $ws.Range("D$r").Value = 0
    $ws.Range("F$r").Value = 0
    $ws.Range("G$r").Value = 15
}

# Update the active selection left behind on the sheet.
$ws.Range("O12").Select()
